# Updates the cryptos worksheet with refreshed prices / 1h volume figures
# (mirrors the scheduled "Updated cryptos list ... with GitHub Actions" commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D/E are plain text in this sheet (prices use "." as a thousands
# separator and volumes are padded percentage strings), so force the cells
# that would otherwise be auto-parsed as numbers to stay Text before writing.
$numericLookingPriceCells = @("D5","D6","D7","D14","D16","D21","D22","D24","D26","D29","D30","D31","D38","D41","D42","D44","D45","D46","D48","D49","D51")
foreach ($addr in $numericLookingPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "61.156.09"
$ws.Range("E2").Value = "  -0.43%  "
# Row 3
$ws.Range("D3").Value = "2.376.29"
$ws.Range("E3").Value = "  -0.94%  "
# Row 4
$ws.Range("E4").Value = "  +0.08%  "
# Row 5
$ws.Range("D5").Value = "548.96"
$ws.Range("E5").Value = "  -0.25%  "
# Row 6
$ws.Range("D6").Value = "137.99"
$ws.Range("E6").Value = "  -3.20%  "
# Row 7
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.01%  "
# Row 8
$ws.Range("E8").Value = "  -2.48%  "
# Row 9
$ws.Range("D9").Value = "2.377.37"
$ws.Range("E9").Value = "  -0.83%  "
# Row 10
$ws.Range("E10").Value = "  +1.22%  "
# Row 11
$ws.Range("E11").Value = "  +1.36%  "
# Row 12
$ws.Range("E12").Value = "  +1.20%  "
# Row 13
$ws.Range("E13").Value = "  -0.18%  "
# Row 14
$ws.Range("D14").Value = "25.07"
$ws.Range("E14").Value = "  -1.89%  "
# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.790.77"
$ws.Range("E15").Value = "  -1.38%  "
# Row 16
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0000166"
$ws.Range("E16").Value = "  -0.35%  "
# Row 17
$ws.Range("D17").Value = "61.103.60"
$ws.Range("E17").Value = "  +0.18%  "
# Row 18
$ws.Range("D18").Value = "2.389.21"
$ws.Range("E18").Value = "  -0.43%  "
# Row 19
$ws.Range("E19").Value = "  -0.23%  "
# Row 20
$ws.Range("E20").Value = "  -0.25%  "
# Row 21
$ws.Range("D21").Value = "320.88"
$ws.Range("E21").Value = "  +0.37%  "
# Row 22
$ws.Range("D22").Value = "6.70"
$ws.Range("E22").Value = "  -0.58%  "
# Row 24
$ws.Range("D24").Value = "64.24"
$ws.Range("E24").Value = "  +0.71%  "
# Row 25
$ws.Range("E25").Value = "  -12.56%  "
# Row 26
$ws.Range("D26").Value = "8.53"
$ws.Range("E26").Value = "  +3.81%  "
# Row 27
$ws.Range("E27").Value = "  +0.00%  "
# Row 28
$ws.Range("D28").Value = "2.488.01"
$ws.Range("E28").Value = "  -1.07%  "
# Row 29
$ws.Range("D29").Value = "8.15"
$ws.Range("E29").Value = "  +0.02%  "
# Row 30
$ws.Range("D30").Value = "509.28"
$ws.Range("E30").Value = "  -5.39%  "
# Row 31
$ws.Range("B31").Value = "Kaspa"
$ws.Range("C31").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D31").Value = "0.150"
$ws.Range("E31").Value = "  +2.59%  "
# Row 32
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").Value = "0.0₃0886"
$ws.Range("E32").Value = "  -6.59%  "
# Row 33
$ws.Range("E33").Value = "  -4.57%  "
# Row 34
$ws.Range("E34").Value = "  -1.34%  "
# Row 35
$ws.Range("E35").Value = "  -4.71%  "
# Row 36
$ws.Range("E36").Value = "  +0.04%  "
# Row 37
$ws.Range("E37").Value = "  -1.06%  "
# Row 38
$ws.Range("D38").Value = "1.88"
$ws.Range("E38").Value = "  +0.29%  "
# Row 39
$ws.Range("E39").Value = "  +0.17%  "
# Row 40
$ws.Range("E40").Value = "  -4.12%  "
# Row 41
$ws.Range("D41").Value = "18.61"
$ws.Range("E41").Value = "  +2.40%  "
# Row 42
$ws.Range("D42").Value = "146.01"
$ws.Range("E42").Value = "  +5.03%  "
# Row 43
$ws.Range("E43").Value = "  -0.05%  "
# Row 44
$ws.Range("D44").Value = "41.61"
$ws.Range("E44").Value = "  +3.38%  "
# Row 45
$ws.Range("D45").Value = "148.33"
$ws.Range("E45").Value = "  +4.46%  "
# Row 46
$ws.Range("D46").Value = "3.60"
$ws.Range("E46").Value = "  -1.10%  "
# Row 47
$ws.Range("E47").Value = "  -6.63%  "
# Row 48
$ws.Range("D48").Value = "0.0522"
$ws.Range("E48").Value = "  -0.16%  "
# Row 49
$ws.Range("D49").Value = "19.36"
$ws.Range("E49").Value = "  -5.17%  "
# Row 50
$ws.Range("E50").Value = "  -0.72%  "
# Row 51
$ws.Range("D51").Value = "0.0912"
$ws.Range("E51").Value = "  +0.00%  "
